$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.024.49'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.305.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.98'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.32%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.69'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.35%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.30'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.36%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.88'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +8.14%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.664.64'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.279.40'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.806'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.963.78'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.62'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.76%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.16'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.37'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.88'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.85'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.02'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.13'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.80'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.96'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.86'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0696'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.82'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.75'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.14%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.971.96'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.69'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.78'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.85'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.531.92'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.80'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -7.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.59'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.14%  '
